$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new backlog rows (Id 28 and 29) right after the current row 16 ---
$ws.Rows.Item(17).Resize(2).Insert()

# Row 17: Id 28 - "All warning messages needs to be handled"
$ws.Cells.Item(17, 1).Value = 28
$ws.Cells.Item(17, 2).Value = "All warning messages needs to be handled"
$ws.Cells.Item(17, 3).Value = "Medium"
$ws.Cells.Item(17, 4).Value = "Not sprint ready"
$ws.Cells.Item(17, 5).Value = "-"
$ws.Cells.Item(17, 6).Value = 5
$ws.Rows.Item(17).RowHeight = 28.8

# Row 18: Id 29 - "Design buttons to follow new themes and styles"
$ws.Cells.Item(18, 1).Value = 29
$ws.Cells.Item(18, 2).Value = "Design buttons to follow new themes and styles"
$ws.Cells.Item(18, 3).Value = "High"
$ws.Cells.Item(18, 4).Value = "Not sprint ready"
$ws.Cells.Item(18, 5).Value = "Create vector images all color combinations"
$ws.Cells.Item(18, 6).Value = 13
$ws.Rows.Item(18).RowHeight = 28.8

# --- Renumber the trailing placeholder rows that were pushed down by the insert ---
# (previously Id 28/29 blank placeholder rows, now Id 30/31; the final boxed row becomes Id 32)
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Style = "Normal"
$ws.Cells.Item(33, 1).Value = 32
